# Add a new "2023" column (N) to the CITES permits table, mirroring the
# existing "2022" column (M): same header style, same bottom-border style
# on the spacer row, and the new permit count in the data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (thin bottom-bordered spacer row above the header) becomes a touch
# taller with an explicit custom height.
$ws.Rows("3").RowHeight = 13.5

# Copy column M's formatting (styles/borders for rows 3-5) into column N so
# the new cells (N3/N4/N5) pick up the same styles as M3/M4/M5 without
# hand-picking style indices.
$ws.Range("M3:M5").Copy($ws.Range("N3:N5"))

# N3 is just a styled spacer cell - keep it empty (copy brought the style,
# clear the copied value).
$ws.Range("N3").Value = ""

# N4: the new year header.
$ws.Range("N4").Value = 2023

# N5: the new permit count for 2023.
$ws.Range("N5").Value = 553
